$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "user story" header row (A44), styled like rows 14/23/29 (wrap text, taller row)
$ws.Range("A44").Value = "Korisniku blokiranog novcanika dozvoljene su samo operacije upita stanja I pregleda transakcija"
$ws.Range("A44").WrapText = $true
$ws.Rows.Item(44).RowHeight = 30

# New task rows (45-47)
$ws.Range("B45").Value = "Prosirivanje Wallet entiteta IsBlocked property-jem I Block/Unblock metodama"
$ws.Range("C45").Value = 5

$ws.Range("B46").Value = "Dodavanje provere da li je Wallet blokiran u metode za Deposit/Withdraw/Transfer"
$ws.Range("C46").Value = 5

$ws.Range("B47").Value = "Dodavanje testova za blokirani wallet na Deposit/Withdraw/Transfer test klase"
$ws.Range("C47").Value = 10

# Update selection to reflect the new last data cell
$ws.Range("D48").Select()
